# "prepare to merge into dev"
# - Reorder the three NPC dialogue blocks on the "plot" sheet: the 6001
#   block (rows 2-4) moves down to the bottom, below the 105001 block, so
#   the sheet now reads 110001, 105001, 6001 (it used to read 6001,
#   110001, 105001).
# - Move the active tab / selection from "plot" to "gossip".

$wb  = $excel.ActiveWorkbook
$plot   = $wb.Worksheets.Item("plot")
$gossip = $wb.Worksheets.Item("gossip")

# --- Re-write rows 2-10 of "plot" with the rotated block order ------------
# New row 2-4  <- old row 5-7   (110001 block)
# New row 5-7  <- old row 8-10  (105001 block)
# New row 8-10 <- old row 2-4   (6001 block)

$plot.Range("A2").Value = 110001
$plot.Range("B2").Value = 110001
$plot.Range("C2").Value = 110003
$plot.Range("D2").Value = "这个是110001"
$plot.Range("E2").Value = 1002
$plot.Range("F2").Value = "plot"
$plot.Range("G2").Value = "啊草"
$plot.Range("H2").Value = "money"

$plot.Range("A3").ClearContents()
$plot.Range("B3").Value = 110002
$plot.Range("C3").ClearContents()
$plot.Range("D3").Value = "这个是110002"
$plot.Range("E3").Value = 1002
$plot.Range("F3").ClearContents()
$plot.Range("G3").ClearContents()
$plot.Range("H3").ClearContents()

$plot.Range("A4").ClearContents()
$plot.Range("B4").Value = 110003
$plot.Range("C4").ClearContents()
$plot.Range("D4").Value = "这个是110003"
$plot.Range("E4").Value = 1002
$plot.Range("F4").ClearContents()
$plot.Range("G4").ClearContents()
$plot.Range("H4").ClearContents()

$plot.Range("A5").Value = 105001
$plot.Range("B5").Value = 105001
$plot.Range("C5").Value = 105003
$plot.Range("D5").Value = "族长的秋天"
$plot.Range("E5").Value = 1002
$plot.Range("F5").Value = "plot"
$plot.Range("G5").Value = "啊草"
$plot.Range("H5").ClearContents()

$plot.Range("A6").ClearContents()
$plot.Range("B6").Value = 105002
$plot.Range("C6").ClearContents()
$plot.Range("D6").Value = "恶时辰"
$plot.Range("E6").Value = 1002
$plot.Range("F6").ClearContents()
$plot.Range("G6").ClearContents()
$plot.Range("H6").ClearContents()

$plot.Range("A7").ClearContents()
$plot.Range("B7").Value = 105003
$plot.Range("C7").ClearContents()
$plot.Range("D7").Value = "世界上最美丽的溺水者"
$plot.Range("E7").Value = 1002
$plot.Range("F7").ClearContents()
$plot.Range("G7").ClearContents()
$plot.Range("H7").ClearContents()

$plot.Range("A8").Value = 6001
$plot.Range("B8").Value = 6001
$plot.Range("C8").Value = 6003
$plot.Range("D8").Value = "这个是6001"
$plot.Range("E8").Value = 1002
$plot.Range("F8").Value = "plot"
$plot.Range("G8").Value = "啊草"
$plot.Range("H8").ClearContents()

$plot.Range("A9").ClearContents()
$plot.Range("B9").Value = 6002
$plot.Range("C9").ClearContents()
$plot.Range("D9").Value = "这个是6002"
$plot.Range("E9").Value = 1002
$plot.Range("F9").ClearContents()
$plot.Range("G9").ClearContents()
$plot.Range("H9").ClearContents()

$plot.Range("A10").ClearContents()
$plot.Range("B10").Value = 6003
$plot.Range("C10").ClearContents()
$plot.Range("D10").Value = "这个是6003"
$plot.Range("E10").Value = 1002
$plot.Range("F10").ClearContents()
$plot.Range("G10").ClearContents()
$plot.Range("H10").ClearContents()

# --- View state: selection on "plot" moves, "gossip" becomes the active tab
$plot.Range("E17").Select() | Out-Null
$gossip.Activate() | Out-Null
